$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / account holder info -----------------------------------------
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long numeric-looking card number that must stay TEXT (it would
# otherwise be auto-coerced to a number by the General format). Stage it in
# a scratch cell formatted as Text, then paste only the value back onto B3
# so the destination keeps its original style (s="8") while the stored type
# remains a string.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "2570314725427075"
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)
$scratch.Delete()

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---------------------------------------------------
$ws.Range("D5").Value = "KONTOSTAND AM 19.02.2024"

# --- Transaction rows 6-9 (dates / description / amount updated) ----------
$ws.Range("B6").Value = "23.02."
$ws.Range("C6").Value = "24.02."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 63577997"
$ws.Range("E6").Value = "37,60-"

$ws.Range("B7").Value = "24.02."
$ws.Range("C7").Value = "25.02."
$ws.Range("E7").Value = "25,14-"

$ws.Range("B8").Value = "27.02."
$ws.Range("C8").Value = "28.02."
$ws.Range("D8").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E8").Value = "70,48-"

$ws.Range("B9").Value = "29.02."
$ws.Range("C9").Value = "01.03."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-47431116"
$ws.Range("E9").Value = "54,03-"

# --- Row 10 was previously blank; it now holds a new transaction ----------
$ws.Range("B10").Value = "03.03."
$ws.Range("C10").Value = "04.03."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 77479917"

# E10 switches from the blank-row style (s=12) to the filled-row style used
# by E6:E9 (s=17). Pull the format from E9 before writing the value so the
# cell picks up the already-existing style instead of minting a new one.
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "86,84-"

# --- Closing balance / next statement date ---------------------------------
$ws.Range("D12").Value = "KONTOSTAND AM 06.03.2024"
$ws.Range("E12").Value = "274,09-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 15.03.2024"
